$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-5 hold one record each. The records got cyclically rotated across
# rows (every column except B, "Taxonsorteringsordning", moved together as
# a unit), while column B received a brand-new value on every row:
#
#   new row 2 <- old row 3 data   (B becomes 89557)
#   new row 3 <- old row 5 data   (B becomes 77636)
#   new row 4 <- old row 2 data   (B becomes 56446)
#   new row 5 <- old row 4 data   (B becomes 90800)
#
# Capture the "old" values for every relevant column on every row first, so
# the subsequent writes (which happen row by row, top to bottom) never read
# already-overwritten data.

$cols = @('A','B','D','E','F','G','H','Q','R','S','Z','AB','AW','AX')

$old = @{}
foreach ($r in 2..5) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $old[$r] = $rowVals
}

$newB = @{ 2 = 89557; 3 = 77636; 4 = 56446; 5 = 90800 }
$sourceRow = @{ 2 = 3; 3 = 5; 4 = 2; 5 = 4 }

foreach ($r in 2..5) {
    $src = $old[$sourceRow[$r]]

    $ws.Range("A$r").Value = $src['A']
    $ws.Range("B$r").Value = $newB[$r]
    $ws.Range("D$r").Value = $src['D']
    $ws.Range("E$r").Value = $src['E']
    $ws.Range("F$r").Value = $src['F']
    $ws.Range("G$r").Value = $src['G']
    $ws.Range("H$r").Value = $src['H']
    $ws.Range("Q$r").Value = $src['Q']
    $ws.Range("R$r").Value = $src['R']
    $ws.Range("S$r").Value = $src['S']
    $ws.Range("Z$r").Value = $src['Z']
    $ws.Range("AB$r").Value = $src['AB']
    $ws.Range("AW$r").Value = $src['AW']
    $ws.Range("AX$r").Value = $src['AX']
}
